# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the refreshed data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 10535
$ws.Range("F6").Value = 822
$ws.Range("F9").Value = 339
$ws.Range("F10").Value = 1103
$ws.Range("F15").Value = 1802
$ws.Range("F19").Value = 534
$ws.Range("F20").Value = 768
$ws.Range("F21").Value = 869
$ws.Range("F25").Value = 596
$ws.Range("F26").Value = 609
$ws.Range("F27").Value = 105
$ws.Range("F29").Value = 994
$ws.Range("F31").Value = 483
$ws.Range("F35").Value = 540
$ws.Range("F36").Value = 1636
$ws.Range("F37").Value = 361
$ws.Range("F39").Value = 1389
$ws.Range("F41").Value = 114

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 181

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2145
$ws.Range("F4").Value = 522

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2145
$ws.Range("F5").Value = 10535
$ws.Range("F8").Value = 522
$ws.Range("F10").Value = 822
$ws.Range("F11").Value = 181
$ws.Range("F12").Value = 339
$ws.Range("F13").Value = 1103
$ws.Range("F18").Value = 1802
$ws.Range("F22").Value = 534
$ws.Range("F23").Value = 768
$ws.Range("F24").Value = 869
$ws.Range("F28").Value = 596
$ws.Range("F31").Value = 609
$ws.Range("F32").Value = 105
$ws.Range("F34").Value = 994
$ws.Range("F37").Value = 483
$ws.Range("F40").Value = 361
$ws.Range("F41").Value = 1389
$ws.Range("F43").Value = 114

$wb.Save()
